$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert a blank row at position 13, shifting old rows 13-23 to 14-24.
$ws.Rows.Item(13).Insert()

# Step 2: the inserted row 13 only carried over a stray styled A13 cell; remove it
# completely (target layout has no A13 -- only B13/C13).
$ws.Cells.Item(13,1).Clear()

# Step 3: (re)write every cell value, row-major / column A->B->C, matching the
# target order so new shared-strings are appended in the same sequence as the
# authored workbook.
$ws.Cells.Item(1,2).Value = "Ementa atual:"
$ws.Cells.Item(1,3).Value = "Ementa modificada (dados modificados em vermelho):"

$ws.Cells.Item(2,2).Value = "LOQ4072"
$ws.Cells.Item(2,3).Value = "LOQ4072"

$ws.Cells.Item(3,1).Value = "Nome:"
$ws.Cells.Item(3,2).Value = " Processos de Separação por Membranas"
$ws.Cells.Item(3,3).Value = " Processos de Separação por Membranas"

$ws.Cells.Item(4,1).Value = "Name:"
$ws.Cells.Item(4,2).Value = "Membrane Separation Processes"
$ws.Cells.Item(4,3).Value = "Membrane Separation Processes"

$ws.Cells.Item(5,1).Value = "Créditos-aula:"
$ws.Cells.Item(5,2).Value = "4"
$ws.Cells.Item(5,3).Value = "4"

$ws.Cells.Item(6,1).Value = "Créditos-trabalho"
$ws.Cells.Item(6,2).Value = "0"
$ws.Cells.Item(6,3).Value = "0"

$ws.Cells.Item(7,1).Value = "Carga horária:"
$ws.Cells.Item(7,2).Value = "60 h"
$ws.Cells.Item(7,3).Value = "60 h"

$ws.Cells.Item(8,1).Value = "Ativação:"
$ws.Cells.Item(8,2).Value = "01/01/2018"
$ws.Cells.Item(8,3).Value = "01/01/2018"

$ws.Cells.Item(9,1).Value = "Semestre ideal:"
$ws.Cells.Item(9,2).Value = "EQD-10,EQN-12"
$ws.Cells.Item(9,3).Value = "EQD-10,EQN-12"

$ws.Cells.Item(10,1).Value = "Objetivos:"
$ws.Cells.Item(10,2).Value = "Propiciar os conhecimentos básicos dos Princípios da Separação por Membranas. Transmitir conhecimentos gerais que permitam entender a tecnologia envolvida nos diferentes tipos de Processos de Separação por Membranas e conhecer os materiais mais utilizados na fabricação de membranas e seus respectivos processos de fabricação. Estudar os componentes fenomenológicos envolvidos nos mecanismos de transporte através das membranas, e identificar as principais vantagens, desvantagens e aplicações deste tipo de processo de separação."
$ws.Cells.Item(10,3).Value = "Propiciar os conhecimentos básicos dos Princípios da Separação por Membranas. Transmitir conhecimentos gerais que permitam entender a tecnologia envolvida nos diferentes tipos de Processos de Separação por Membranas e conhecer os materiais mais utilizados na fabricação de membranas e seus respectivos processos de fabricação. Estudar os componentes fenomenológicos envolvidos nos mecanismos de transporte através das membranas, e identificar as principais vantagens, desvantagens e aplicações deste tipo de processo de separação."

$ws.Cells.Item(11,1).Value = "Objectives:"
$ws.Cells.Item(11,2).Value = "Provide basic knowledge of the Principles of Membrane Separation. Provide general information to understand the technology involved in the different types of Membrane Separation Processes and to know the materials most used in the manufacture of membranes and their respective manufacturing processes. To study the phenomenological components involved in the mechanisms of transport through the membranes, and to identify the main advantages, disadvantages and applications of this type of separation process."
$ws.Cells.Item(11,3).Value = "Provide basic knowledge of the Principles of Membrane Separation. Provide general information to understand the technology involved in the different types of Membrane Separation Processes and to know the materials most used in the manufacture of membranes and their respective manufacturing processes. To study the phenomenological components involved in the mechanisms of transport through the membranes, and to identify the main advantages, disadvantages and applications of this type of separation process."

$ws.Cells.Item(12,1).Value = "Docentes responsáveis:"

$ws.Cells.Item(13,2).Value = "787307 - Luis Fernando Figueiredo Faria"
$ws.Cells.Item(13,3).Value = "787307 - Luis Fernando Figueiredo Faria"

$ws.Cells.Item(14,1).Value = "Programa resumido:"
$ws.Cells.Item(14,2).Value = "Classificação dos processos com membranas e suas aplicações. Técnicas de preparo dos diferentes tipos de membranas poliméricas. Tipos de módulos e suas principais características. Fundamentos teóricos da síntese de membranas pela inversão de fases; influência das variáveis de síntese nas características de transporte das membranas. Apresentação dos diferentes tipos de processos com membranas. Aplicações. Projeto para uma aplicação específica."
$ws.Cells.Item(14,3).Value = "Classificação dos processos com membranas e suas aplicações. Técnicas de preparo dos diferentes tipos de membranas poliméricas. Tipos de módulos e suas principais características. Fundamentos teóricos da síntese de membranas pela inversão de fases; influência das variáveis de síntese nas características de transporte das membranas. Apresentação dos diferentes tipos de processos com membranas. Aplicações. Projeto para uma aplicação específica."

$ws.Cells.Item(15,1).Value = "Short syllabus:"
$ws.Cells.Item(15,2).Value = "Classification of membrane processes and their applications. Technical preparation of different types of polymeric membranes. Types of modules and its main features. Theoretical fundamentals of the membranes synthesis  by the inversion of phases; Influence of the synthesis variables on the transport characteristics of the membranes. Presentation of the different types of membrane processes. Applications. Design for a specific application."
$ws.Cells.Item(15,3).Value = "Classification of membrane processes and their applications. Technical preparation of different types of polymeric membranes. Types of modules and its main features. Theoretical fundamentals of the membranes synthesis  by the inversion of phases; Influence of the synthesis variables on the transport characteristics of the membranes. Presentation of the different types of membrane processes. Applications. Design for a specific application."

$ws.Cells.Item(16,1).Value = "Programa:"
$ws.Cells.Item(16,2).Value = "1. Introdução: Histórico e definição de processos de separação com membranas: comparação com processos clássicos de separação. Processos com membranas mais utilizados: classificação e aplicações. 2. Membranas: Definição; materiais empregados; classificação quanto à estrutura e quanto ao tipo de separação a que se destina; membranas microporosas: técnicas de fabricação e características. Síntese de membranas pela técnica da inversão de fases:aspectos termodinâmicos e cinéticos.3. Fundamentos dos Processos de Separação por Membranas: mecanismos de transporte em membranas densas e microporosas.4. Membranas e Módulos Comerciais: geometria de membranas; técnicas de fabricação de membranas planas, membranas do tipo fibra oca e tubulares. Módulos: tipo placa, espiral, tubular e fibra oca.5. Processos Comerciais de Separação com Membranas. Principais características aplicações."
$ws.Cells.Item(16,3).Value = "1. Introdução: Histórico e definição de processos de separação com membranas: comparação com processos clássicos de separação. Processos com membranas mais utilizados: classificação e aplicações. 2. Membranas: Definição; materiais empregados; classificação quanto à estrutura e quanto ao tipo de separação a que se destina; membranas microporosas: técnicas de fabricação e características. Síntese de membranas pela técnica da inversão de fases:aspectos termodinâmicos e cinéticos.3. Fundamentos dos Processos de Separação por Membranas: mecanismos de transporte em membranas densas e microporosas.4. Membranas e Módulos Comerciais: geometria de membranas; técnicas de fabricação de membranas planas, membranas do tipo fibra oca e tubulares. Módulos: tipo placa, espiral, tubular e fibra oca.5. Processos Comerciais de Separação com Membranas. Principais características aplicações."

$ws.Cells.Item(17,1).Value = "Syllabus:"
$ws.Cells.Item(17,2).Value = "1. Introduction: History and definition of membranes separation processes: comparison with classical separation processes. Membrane processes: classification and applications. 2. Membranes: Definition; materials used; rating on the structure and the type of proposed separation; Microporous membranes: manufacturing techniques and characteristics. Synthesis of membranes by the phase inversion technique: thermodynamic and kinetic aspects. 3. Fundamentals of Membrane Separation Processes: transport mechanisms in dense and microporous membranes.4. Membranes and Commercial Modules: membrane geometry; Techniques for manufacturing flat membranes, hollow fiber and tubular membranes. Modules: plate type, spiral, tubular and hollow fiber.5. Commercial Membrane Separation Processes. Main applications."
$ws.Cells.Item(17,3).Value = "1. Introduction: History and definition of membranes separation processes: comparison with classical separation processes. Membrane processes: classification and applications. 2. Membranes: Definition; materials used; rating on the structure and the type of proposed separation; Microporous membranes: manufacturing techniques and characteristics. Synthesis of membranes by the phase inversion technique: thermodynamic and kinetic aspects. 3. Fundamentals of Membrane Separation Processes: transport mechanisms in dense and microporous membranes.4. Membranes and Commercial Modules: membrane geometry; Techniques for manufacturing flat membranes, hollow fiber and tubular membranes. Modules: plate type, spiral, tubular and hollow fiber.5. Commercial Membrane Separation Processes. Main applications."

$ws.Cells.Item(18,1).Value = "Avaliação:"

$ws.Cells.Item(19,1).Value = "Método:"
$ws.Cells.Item(19,2).Value = "-Provas escritas; -participação e conteúdo de trabalho e seminário;"
$ws.Cells.Item(19,3).Value = "-Provas escritas; -participação e conteúdo de trabalho e seminário;"

$ws.Cells.Item(20,1).Value = "Critério:"
$ws.Cells.Item(20,2).Value = "Média Final = (Prova1 + Prova2 + Nota de Trabalho) / 3Média final mínima de aprovação = 5,0"
$ws.Cells.Item(20,3).Value = "Média Final = (Prova1 + Prova2 + Nota de Trabalho) / 3Média final mínima de aprovação = 5,0"

$ws.Cells.Item(21,1).Value = "Norma de recuperação:"
$ws.Cells.Item(21,2).Value = "(Prova escrita + Média Final)/2         Nota Final mínima para aprovação= 5,0"
$ws.Cells.Item(21,3).Value = "(Prova escrita + Média Final)/2         Nota Final mínima para aprovação= 5,0"

$ws.Cells.Item(22,1).Value = "Bibliografia:"
$ws.Cells.Item(22,2).Value = "1 - HABERT, A. C.; BORGES, C. P.; NÓBREGA, R. Processos de separação por membranas. Rio de Janeiro: E-papers, 2006. 180p.`n2 - BAKER, R. W. Membrane Technology and Applications. 2nd. ed. Chichester : John Wiley & Sons, 2004. `n3 - MULDER, M. Basic Principles of Membrane Technology. Holanda: Klumer Academic Publishers, 1991.`n4 - CHERYAN, M. Ultrafiltration and Microfiltration Handbook. USA: Technomic Publishing Co. Inc, 1998.`n5 - Membrane Handbook. Ed. W.S.W Ho and K.K. Sirkar. New York : Chapman & Hall, 1992.`n6 - RAUTENBACH, R.; ALBRECHT, R. Membrane Processes / Ed. Antony Rowe Ltd. Wiltshire, Great Britain, 1994."
$ws.Cells.Item(22,3).Value = "1 - HABERT, A. C.; BORGES, C. P.; NÓBREGA, R. Processos de separação por membranas. Rio de Janeiro: E-papers, 2006. 180p.`n2 - BAKER, R. W. Membrane Technology and Applications. 2nd. ed. Chichester : John Wiley & Sons, 2004. `n3 - MULDER, M. Basic Principles of Membrane Technology. Holanda: Klumer Academic Publishers, 1991.`n4 - CHERYAN, M. Ultrafiltration and Microfiltration Handbook. USA: Technomic Publishing Co. Inc, 1998.`n5 - Membrane Handbook. Ed. W.S.W Ho and K.K. Sirkar. New York : Chapman & Hall, 1992.`n6 - RAUTENBACH, R.; ALBRECHT, R. Membrane Processes / Ed. Antony Rowe Ltd. Wiltshire, Great Britain, 1994."

$ws.Cells.Item(23,1).Value = "Requisitos:"

$ws.Cells.Item(24,2).Value = "LOQ4085 -  Operações Unitárias I  (Requisito fraco)`n"
$ws.Cells.Item(24,3).Value = "LOQ4085 -  Operações Unitárias I  (Requisito fraco)`n"

# Step 4: row 13 (B/C) inherited column-As bold style from the row insert;
# repaint B13/C13 with the correct column styles by copying formats from a
# same-column cell elsewhere (reuses existing cellXfs, no new styles created).
$ws.Cells.Item(3,2).Copy()
$ws.Cells.Item(13,2).PasteSpecial(-4122)
$ws.Cells.Item(3,3).Copy()
$ws.Cells.Item(13,3).PasteSpecial(-4122)

Write-Output "edit complete"
